# Particles.js background en login fixes
# Applies the Bestand2 sheet changes:
#  - restyle the "Bestand 1 column N" fill from blue to pink
#  - strip the leftover "Bestand 2" text that had been concatenated onto
#    the "Bestand 1 column N" strings in column C
#  - add a new column AB ("Bestand 2 column N" helper data) with its own
#    width, and move the selection/scroll position over to it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up column C: drop the accidental "Bestand 2" prefix ---------
$fixedC = @{
    "C4"  = "Bestand 1 column 4"
    "C6"  = "Bestand 1 column 6"
    "C7"  = "Bestand 1 column 2"
    "C9"  = "Bestand 1 column 8"
    "C10" = "Bestand 1 column 11"
    "C11" = "Bestand 1 column 10"
    "C12" = "Bestand 1 column 9"
}
foreach ($addr in $fixedC.Keys) {
    $ws.Range($addr).Value = $fixedC[$addr]
}

# --- 2. Recolor those same cells: blue (0069FF) -> pink (FF9F9F) -------
foreach ($addr in $fixedC.Keys) {
    $ws.Range($addr).Interior.Color = 10461183   # RGB(255, 159, 159)
}

# --- 3. Add new column AB (28) with "Bestand 2" data -------------------
$colAB = @{
    1  = "Dit is bestand 2 column 28"
    2  = "Bestand 2 column 2"
    3  = "Bestand 2 column 3"
    4  = "Bestand 2 column 4"
    5  = "Bestand 2 column 5"
    6  = "Bestand 2 column 6"
    7  = "Bestand 2 column 7"
    8  = "Bestand 2 column 8"
    9  = "Bestand 2 column 9"
    10 = "Bestand 2 column 10"
    11 = "Bestand 2 column 11"
    12 = "Bestand 2 column 12"
}
foreach ($row in $colAB.Keys) {
    $ws.Cells.Item($row, 28).Value = $colAB[$row]
}

# Widen column AB to match the other data columns.
$ws.Columns.Item(28).ColumnWidth = 31.5

# --- 4. Move the view over to the new data: scroll + select AB5 --------
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("AB5").Select()
